# Include FLD_Transmittals_OverDue_New_IssuedForReview for execution.
#
# 1) Rows 12-14 ("ActionRequired_CaC_*" test cases) move from Sprint3 to Sprint5.
# 2) A new test case row is appended (row 16):
#    FLD_Transmittals_OverDue_New_IssuedForReview /
#    "Verifies the transmittal under Actions Overdue menu", scheduled for Sprint5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-schedule the three "Close and Cancel" cases from Sprint3 to Sprint5 ---
$ws.Range("F12").Value = "Sprint5"
$ws.Range("F13").Value = "Sprint5"
$ws.Range("F14").Value = "Sprint5"

# --- 2. Append a new row for the OverDue/IssuedForReview test case ---
# Copy the formatting of the last existing row (15) down into the new row (16)
# so borders / wrap-text / alignment stay consistent with the rest of the table.
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "FLD_Transmittals_OverDue_New_IssuedForReview"
$ws.Range("B16").Value = "Verifies the transmittal under Actions Overdue menu"
$ws.Range("C16").Value = "Y"
$ws.Range("D16").Value = "Y"
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = "Sprint5"

# --- 3. Extend the two data-validation list ranges so they cover the new row ---
$ws.Range("C2:D16").Validation.Delete()
$ws.Range("C2:D16").Validation.Add(3, 1, 3, '"Y,N"')

$ws.Range("F2:F16").Validation.Delete()
$ws.Range("F2:F16").Validation.Add(3, 1, 3, '"Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10"')
